$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 21.47 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 1182.86 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 1506 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 10188 -> 11811
$t.Cell(4, 1).Range.Text = "11811"

# Row 7: 0.08368 -> 0.08601
$t.Cell(7, 1).Range.Text = "0.08601"

# Row 8: 0.06933 -> 0.07376
$t.Cell(8, 1).Range.Text = "0.07376"

# Row 12: 976.44867 -> 1182.85856
$t.Cell(12, 1).Range.Text = "1182.85856"

# Row 44: collapse the multi-run tab-separated values down to a single value 21.47
$t.Cell(44, 1).Range.Text = "21.47"

# Row 45: collapse the multi-run tab-separated values down to a single value 1182.86
$t.Cell(45, 1).Range.Text = "1182.86"

# Row 46: collapse the multi-run tab-separated values down to a single value 1506
$t.Cell(46, 1).Range.Text = "1506"
